$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 366; all existing rows from 366 downward
# (including the former last row 388) shift down by one, to 389.
$ws.Rows.Item(366).Insert()

# Populate the newly inserted row 366 with the new weekly record.
$ws.Range("A366").Value = 1
$ws.Range("B366").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C366").Value = "Arica y Parinacota"
$ws.Range("D366").Value = 44931
$ws.Range("E366").Value = 15
$ws.Range("F366").Value = 100114013
$ws.Range("G366").Value = "Zanahoria"
$ws.Range("H366").Value = "Sin especificar"
$ws.Range("I366").Value = "Primera"
$ws.Range("J366").Value = 35
$ws.Range("K366").Value = 26000
$ws.Range("L366").Value = 27000
$ws.Range("M366").Value = 26429
$ws.Range("N366").Value = "$/saco 25 kilos"
$ws.Range("O366").Value = "Valle de Camiña"
$ws.Range("P366").Value = 1057
$ws.Range("Q366").Value = 25
$ws.Range("R366").Value = "Hortaliza"
